# Region IX_HEALTH.xlsx update:
#  - Header text casing tweaks (I1, L1)
#  - Remove placeholder "-" values from I2:I9 and L2:L9 (now blank cells)
#  - Insert five new tracking-count columns (AA:AE) before the existing
#    "Status as of ..." column, which shifts it from AA to AF and extends
#    the dimension / data validation accordingly
#  - Populate the new header cells

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header text changes -------------------------------------------------
$ws.Range("I1").Value = "TOTAL PHYSICAL TARGET"
$ws.Range("L1").Value = "BATCH"

# --- Clear stale "-" placeholder cells in I and L for rows 2-9 ----------
$ws.Range("I2:I9").ClearContents()
$ws.Range("L2:L9").ClearContents()

# --- Insert 5 new columns before the old AA column ("Status as of ...") -
# This shifts the existing AA column (and its data validation / dimension)
# to AF, and the newly inserted columns inherit the header formatting
# from the column immediately to their left (Z), matching style s="1".
$ws.Range("AA1:AE1").EntireColumn.Insert()

# --- Populate the new header cells (AA1:AE1) -----------------------------
$ws.Range("AA1").Value = "No. of Sites Reverted"
$ws.Range("AB1").Value = "No. of Sites Not yet started"
$ws.Range("AC1").Value = "No. of Sites Under Procurement"
$ws.Range("AD1").Value = "No. of Sites On Going"
$ws.Range("AE1").Value = "No. of Sites Completed"
